$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data row for account 004452912 / Bruno / 200000 (row 2, right
# after the header row). Deleting the entire row shifts everything below it
# up by one, matching the rest of the table staying intact.
$ws.Rows.Item(2).Delete()
